$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
$c1 = $tcs.Colors(1)
$c1.RGB = 255
